# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures on the active sheet to the latest scraped snapshot.
#
# Several "Price" values look numeric (e.g. "0.0791", "212.56") but must stay
# plain text, matching how the sheet already stores them (no thousands/decimal
# coercion, leading zeros preserved, etc.). Assigning such a string directly
# via .Value lets Excel auto-detect it as a number, so for those cells we
# write it as a text-forced ("quote-prefixed") entry and then immediately
# reapply the "Normal" cell style so no stray number-format is left behind on
# the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '25.938.67'; ForceText = $false },
    @{ Cell = "E2"; Value = '  -0.43%  '; ForceText = $false },
    @{ Cell = "D3"; Value = '1.621.15'; ForceText = $false },
    @{ Cell = "E3"; Value = '  -1.09%  '; ForceText = $false },
    @{ Cell = "E4"; Value = '  -0.37%  '; ForceText = $false },
    @{ Cell = "D5"; Value = '212.56'; ForceText = $true },
    @{ Cell = "E5"; Value = '  -0.99%  '; ForceText = $false },
    @{ Cell = "E7"; Value = '  -0.37%  '; ForceText = $false },
    @{ Cell = "E8"; Value = '  -0.89%  '; ForceText = $false },
    @{ Cell = "E9"; Value = '  -1.51%  '; ForceText = $false },
    @{ Cell = "D10"; Value = '18.41'; ForceText = $true },
    @{ Cell = "E10"; Value = '  -1.32%  '; ForceText = $false },
    @{ Cell = "D11"; Value = '0.0791'; ForceText = $true },
    @{ Cell = "E11"; Value = '  -0.39%  '; ForceText = $false },
    @{ Cell = "D12"; Value = '1.846.88'; ForceText = $false },
    @{ Cell = "E12"; Value = '  -1.09%  '; ForceText = $false },
    @{ Cell = "D13"; Value = '1.632.01'; ForceText = $false },
    @{ Cell = "E13"; Value = '  -3.21%  '; ForceText = $false },
    @{ Cell = "E14"; Value = '  -1.64%  '; ForceText = $false },
    @{ Cell = "D15"; Value = '0.524'; ForceText = $true },
    @{ Cell = "E15"; Value = '  -1.37%  '; ForceText = $false },
    @{ Cell = "D16"; Value = '25.971.25'; ForceText = $false },
    @{ Cell = "E16"; Value = '  -0.38%  '; ForceText = $false },
    @{ Cell = "D17"; Value = '61.67'; ForceText = $true },
    @{ Cell = "E17"; Value = '  -1.08%  '; ForceText = $false },
    @{ Cell = "E18"; Value = '  -1.42%  '; ForceText = $false },
    @{ Cell = "E19"; Value = '  -0.37%  '; ForceText = $false },
    @{ Cell = "D20"; Value = '191.81'; ForceText = $true },
    @{ Cell = "E20"; Value = '  +0.37%  '; ForceText = $false },
    @{ Cell = "D21"; Value = '4.25'; ForceText = $true },
    @{ Cell = "E21"; Value = '  -0.48%  '; ForceText = $false },
    @{ Cell = "E22"; Value = '  -0.85%  '; ForceText = $false },
    @{ Cell = "D23"; Value = '6.02'; ForceText = $true },
    @{ Cell = "E23"; Value = '  -2.16%  '; ForceText = $false },
    @{ Cell = "E24"; Value = '  +1.07%  '; ForceText = $false },
    @{ Cell = "D25"; Value = '143.74'; ForceText = $true },
    @{ Cell = "E25"; Value = '  -0.34%  '; ForceText = $false },
    @{ Cell = "E26"; Value = '  -0.46%  '; ForceText = $false },
    @{ Cell = "D27"; Value = '1.71'; ForceText = $true },
    @{ Cell = "E27"; Value = '  -3.00%  '; ForceText = $false },
    @{ Cell = "E28"; Value = '  -2.20%  '; ForceText = $false },
    @{ Cell = "E29"; Value = '  -0.09%  '; ForceText = $false },
    @{ Cell = "E30"; Value = '  -1.11%  '; ForceText = $false },
    @{ Cell = "D31"; Value = '0.0479'; ForceText = $true },
    @{ Cell = "E31"; Value = '  -1.58%  '; ForceText = $false },
    @{ Cell = "E32"; Value = '  -1.59%  '; ForceText = $false },
    @{ Cell = "E33"; Value = '  -2.58%  '; ForceText = $false },
    @{ Cell = "D34"; Value = '1.49'; ForceText = $true },
    @{ Cell = "E34"; Value = '  -1.12%  '; ForceText = $false },
    @{ Cell = "E35"; Value = '  -0.77%  '; ForceText = $false },
    @{ Cell = "D36"; Value = '1.128.30'; ForceText = $false },
    @{ Cell = "E36"; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = "E37"; Value = '  -3.76%  '; ForceText = $false },
    @{ Cell = "E38"; Value = '  -2.07%  '; ForceText = $false },
    @{ Cell = "E39"; Value = '  -2.03%  '; ForceText = $false },
    @{ Cell = "E40"; Value = '  -1.43%  '; ForceText = $false },
    @{ Cell = "D41"; Value = '97.67'; ForceText = $true },
    @{ Cell = "E41"; Value = '  -1.37%  '; ForceText = $false },
    @{ Cell = "D42"; Value = '1.757.91'; ForceText = $false },
    @{ Cell = "E42"; Value = '  -0.88%  '; ForceText = $false },
    @{ Cell = "D43"; Value = '0.758'; ForceText = $true },
    @{ Cell = "E43"; Value = '  -3.59%  '; ForceText = $false },
    @{ Cell = "E44"; Value = '  -4.08%  '; ForceText = $false },
    @{ Cell = "E45"; Value = '  -0.72%  '; ForceText = $false },
    @{ Cell = "E46"; Value = '  +1.33%  '; ForceText = $false },
    @{ Cell = "E47"; Value = '  -2.38%  '; ForceText = $false },
    @{ Cell = "D48"; Value = '0.0517'; ForceText = $true },
    @{ Cell = "E48"; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = "E49"; Value = '  -1.03%  '; ForceText = $false },
    @{ Cell = "D50"; Value = '7.49'; ForceText = $true },
    @{ Cell = "E50"; Value = '  -1.24%  '; ForceText = $false },
    @{ Cell = "E51"; Value = '  -0.22%  '; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).Value = "'" + $u.Value
        $ws.Range($u.Cell).Style = "Normal"
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
